# Add an "Academic Number" column (D) to the registered-ids sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in D1, numeric academic numbers in D2:D3
$ws.Range("D1").Value = "Academic Number"
$ws.Range("D2").Value = 12345
$ws.Range("D3").Value = 77225

# Size column D to fit its content, like the existing bestFit column A
$ws.Columns.Item(4).ColumnWidth = 17.5703125

# Mirror the author's final selection on the new header cell
$ws.Range("D1").Select() | Out-Null
